$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brand")

# Remove the existing "Nike" brand row (row 2: ID=1, Name=Nike) by clearing
# its contents so the row disappears from the sheet (rows below keep their
# original row numbers, matching how row 4 is already missing).
$ws.Range("A2:B2").ClearContents()

# Append the brand catalogue with the new entries: Nike (again, now as
# id 6) and the new "emo" brand (id 7).
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Nike"
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "emo"
